$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.882.07'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.640.29'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5061'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2580'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06442'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07788'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.289'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '1.641.57'
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('D14').Value = '1.864.66'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5623'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.73%  '
$ws.Range('D16').Value = '0.0₅7627'
$ws.Range('E16').Value = '  -0.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.19'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '25.897.40'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.329'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.909'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.132'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.777'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1267'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.833'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.244'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04884'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.303'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.229'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.570'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.373'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9044'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.577'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5546'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.128.72'
$ws.Range('E39').Value = '  +0.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01563'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9955'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.547'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8035'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '98.01'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '1.774.54'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  -9.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4276'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.752'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05048'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.33%  '
